$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(237).Insert()

$ws.Cells.Item(237,1).Value = 8
$ws.Cells.Item(237,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(237,3).Value = "Coquimbo"
$ws.Cells.Item(237,4).Value = 44918
$ws.Cells.Item(237,5).Value = 4
$ws.Cells.Item(237,6).Value = 100112037
$ws.Cells.Item(237,7).Value = "Cebollín"
$ws.Cells.Item(237,8).Value = "Sin especificar"
$ws.Cells.Item(237,9).Value = "Primera"
$ws.Cells.Item(237,10).Value = 1200
$ws.Cells.Item(237,11).Value = 1200
$ws.Cells.Item(237,12).Value = 1400
$ws.Cells.Item(237,13).Value = 1300
$ws.Cells.Item(237,14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(237,15).Value = "Provincia del Elquí"
$ws.Cells.Item(237,16).Value = 217
$ws.Cells.Item(237,17).Value = 6
$ws.Cells.Item(237,18).Value = "Hortaliza"
